# Generate Report for handback
#
# Refresh the generated handoff/handback timestamps recorded for the
# "cb4966ff-2498-4ae3-a483-bb23c0f4525d" entry (the last data row, row 4)
# on both the "zh-cn" and "de-de" report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-18 05:37:35"   # Correspond Handoff Datetime
$wsZhCn.Range("G4").Value = "2016-01-18 05:38:28"   # Correspond Handback DateTime

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-18 05:37:47"   # Correspond Handoff Datetime
$wsDeDe.Range("G4").Value = "2016-01-18 05:38:50"   # Correspond Handback DateTime
